# Unit 7 Examples / DynamicArrayOfInt_time-estimates.xlsx
# - Rename "v1.0" -> "v1.0 - 62yrs"
# - Rename "v1.1 ArrayList<Integer>" -> "v1.1 - presumaby around 50yrs"
# - Add a new sheet "v1.2 - 37yrs" (a variant of the v1.1 forecast sheet
#   with an extra two years of data, 8 points instead of 6)
# - Update the lingering selection on the v1.1 sheet

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "v1.0 - 62yrs"
$ws2.Name = "v1.1 - presumaby around 50yrs"

# Build the new "v1.2" sheet by duplicating the v1.1 sheet (same layout,
# headers, styles, and formulas), placed after it, then adjust the data
# that differs.
[void]$ws2.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "v1.2 - 37yrs"

# Column A holds the measured times (ticks) per run; v1.2 adds two more
# data points (rows 8 & 9) and uses different totals throughout.
$ws3.Range("A2").Formula = "=47*2"
$ws3.Range("A3").Formula = "=172+171"
$ws3.Range("A4").Formula = "=641+641"
$ws3.Range("A5").Formula = "=2793+2751"
$ws3.Range("A6").Formula = "=15+11833+11643"
$ws3.Range("A7").Formula = "=16+48711+47524"
$ws3.Range("A8").Formula = "=62+222730+217416"
$ws3.Range("A9").Formula = "=125+1270821+1256945"

# Column B (log2 of A) now extends to rows 8 & 9 to match.
$ws3.Range("B8").Formula = "=LOG(A8,2)"
$ws3.Range("B9").Formula = "=LOG(A9,2)"

# Column E's FORECAST.LINEAR regression now spans the 8 populated points
# (rows 2:9) instead of the original 6 (rows 2:7).
for ($r = 2; $r -le 18; $r++) {
  $ws3.Cells.Item($r, 5).Formula = "=FORECAST.LINEAR(D$r,B`$2:B`$9,D`$2:D`$9)"
}

[void]$ws3.Range("E30").Select()

[void]$ws2.Activate()
[void]$ws2.Range("D30").Select()
